$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.283.26'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.048.19'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.49'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.620'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.99'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.384'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0769'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.76'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.349.82'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.52'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.754'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.80%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.052.48'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.285.22'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.97'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.65'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0822'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.98'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.35'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.49'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.64'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.130'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.13'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.121'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0624'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.57'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.46'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.82'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.25'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0226'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.81%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.488.13'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.27%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.74'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.73%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0951'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.73%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.90'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.45%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.57'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.02'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.88%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.20'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.93'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.237.10'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.45%  '
